# "Loan RBI, Variable Instalments"
#
# The repayment schedule gains a new (currently blank) column so a
# "Variable Instalments" value can be tracked alongside the existing
# Principal/Interest/Fees/Penalties breakdown. Concretely this is a plain
# "insert column" in front of column N ("Late") on the "Repayment schedule"
# sheet - it shifts the old N/O/P columns ("Late"/heading/"Outstanding")
# one slot to the right (-> O/P/Q) and leaves a blank column N behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Remember column M's width so the freshly inserted column N can inherit
# it, same as Excel does automatically when you insert a column.
$mWidth = $ws.Columns("M").ColumnWidth

$ws.Columns("N").Insert() | Out-Null

$ws.Columns("N").ColumnWidth = $mWidth

# Make "Repayment schedule" the active sheet/tab, with the cursor parked
# on S7 (previously "Transactions" was the active tab).
$ws.Activate() | Out-Null
$ws.Range("S7").Select() | Out-Null
